$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the "H 72" record), shifting all subsequent rows up by one.
$ws.Rows.Item(2).Delete()
